$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 251, pushing the existing rows 251..299 down to 252..300.
# This reproduces the observed pattern where every existing record's date /
# price fields shift down one row, and a brand-new weekly record is placed
# at row 251 (the most recent survey date).
$ws.Rows.Item(251).Insert()

# Populate the new row 251 with the latest weekly record.
$ws.Cells.Item(251, 1).Value = 4
$ws.Cells.Item(251, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(251, 3).Value = "Los Lagos"
$ws.Cells.Item(251, 4).Value = 44711
$ws.Cells.Item(251, 5).Value = 10
$ws.Cells.Item(251, 6).Value = 100114014
$ws.Cells.Item(251, 7).Value = "Betarraga"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 500
$ws.Cells.Item(251, 11).Value = 1000
$ws.Cells.Item(251, 12).Value = 1200
$ws.Cells.Item(251, 13).Value = 1100
$ws.Cells.Item(251, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(251, 15).Value = "Región del Maule"
$ws.Cells.Item(251, 16).Value = 220
$ws.Cells.Item(251, 17).Value = 5
$ws.Cells.Item(251, 18).Value = "Hortaliza"
